$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("G2").Value = 2987.520456868959
